$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Actualiza base de datos EC": the "Periodo Mora" column (E16:E22) is being
# refreshed - the list of billing periods is re-entered in chronological
# order (oldest -> newest) instead of the previous reverse-chronological
# order.
$ws.Range("E16").Value = "2308"
$ws.Range("E17").Value = "2309"
$ws.Range("E18").Value = "2310"
$ws.Range("E19").Value = "2311"
$ws.Range("E20").Value = "2312"
$ws.Range("E21").Value = "2401"
$ws.Range("E22").Value = "2402"
